# The three observation rows (20-22) were re-ordered/re-matched: the data
# that used to live in row 21 now belongs to row 20, the data that used to
# live in row 22 now belongs to row 21, and the data that used to live in
# row 20 now belongs to row 22 (a 3-way cyclic rotation). Apply the new
# values cell-by-cell (only the cells that actually change are touched).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 20 (becomes the old row 21's record) ---
$ws.Range("A20").Value = 111559701
$ws.Range("B20").Value = 12249
$ws.Range("D20").Value = "EN"
$ws.Range("E20").Value = 101283
$ws.Range("F20").Value = "Djupsvart brunbagge"
$ws.Range("G20").Value = "Melandrya dubia"
$ws.Range("H20").Value = "(Schaller, 1783)"
$ws.Range("L20").Value = ""
$ws.Range("M20").Value = ""
$ws.Range("Q20").Value = 523950.9321204902
$ws.Range("R20").Value = 6934675.944620069
# Leading apostrophe forces these YYYY-MM-DD looking values to stay text
# (matching the source workbook) instead of being auto-parsed as dates.
$ws.Range("Y20").Value = "'2023-08-17"
$ws.Range("AA20").Value = "'2023-08-17"
$ws.Range("AC20").Value = "Kläckhål med svartfärgade larvgångar på björkhögstubbe med levande fnösktickor. Naturskog norr om Vattensjöarna"

# --- Row 21 (becomes the old row 22's record) ---
$ws.Range("A21").Value = 111560043
$ws.Range("B21").Value = 78578
$ws.Range("D21").Value = "NT"
$ws.Range("E21").Value = 6458
$ws.Range("F21").Value = "Lunglav"
$ws.Range("G21").Value = "Lobaria pulmonaria"
$ws.Range("H21").Value = "(L.) Hoffm."
$ws.Range("L21").ClearContents()
$ws.Range("M21").ClearContents()
$ws.Range("Q21").Value = 523949.236686704
$ws.Range("R21").Value = 6934654.704083432
$ws.Range("Y21").Value = "'2023-08-18"
$ws.Range("AA21").Value = "'2023-08-18"
$ws.Range("AC21").ClearContents()

# --- Row 22 (becomes the old row 20's record) ---
$ws.Range("A22").Value = 111560058
$ws.Range("Q22").Value = 523906.9737172622
$ws.Range("R22").Value = 6934619.326478666
